$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($Row1, $Row2, $FirstCol, $LastCol) {
    for ($col = $FirstCol; $col -le $LastCol; $col++) {
        $c1 = $ws.Cells.Item($Row1, $col)
        $c2 = $ws.Cells.Item($Row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# Columns F (6) through V (22) get swapped between the paired rows.
Swap-RowRange 84 85 6 22
Swap-RowRange 136 137 6 22
Swap-RowRange 140 141 6 22

# Append two brand-new match rows (160, 161) after the former last row (159).
# Copy formatting (styles) from row 159 first, then fill in the values.
$ws.Range("A159:V159").Copy() | Out-Null
$ws.Range("A160:V161").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(160,1).Value = 159
$ws.Cells.Item(160,2).Value = "argentina"
$ws.Cells.Item(160,3).Value = "copa-de-la-liga-profesional"
$ws.Cells.Item(160,4).NumberFormat = "@"
$ws.Cells.Item(160,4).Value = "2023"
$ws.Cells.Item(160,4).Style = "Normal"
$ws.Cells.Item(160,5).Value = 45235.83333333334
$ws.Cells.Item(160,6).Value = "Velez Sarsfield"
$ws.Cells.Item(160,7).Value = 1
$ws.Cells.Item(160,8).Value = "Talleres Cordoba"
$ws.Cells.Item(160,9).Value = 1
$ws.Cells.Item(160,10).Value = 2.47
$ws.Cells.Item(160,11).Value = "30/10/2023 20:12"
$ws.Cells.Item(160,12).Value = 2.22
$ws.Cells.Item(160,13).Value = "05/11/2023 19:57"
$ws.Cells.Item(160,14).Value = 3
$ws.Cells.Item(160,15).Value = "30/10/2023 20:12"
$ws.Cells.Item(160,16).Value = 3.08
$ws.Cells.Item(160,17).Value = "05/11/2023 19:56"
$ws.Cells.Item(160,18).Value = 3.28
$ws.Cells.Item(160,19).Value = "30/10/2023 20:12"
$ws.Cells.Item(160,20).Value = 3.87
$ws.Cells.Item(160,21).Value = "05/11/2023 19:57"
$ws.Cells.Item(160,22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/velez-sarsfield-talleres-cordoba/beG7S4rM/"

$ws.Cells.Item(161,1).Value = 160
$ws.Cells.Item(161,2).Value = "argentina"
$ws.Cells.Item(161,3).Value = "copa-de-la-liga-profesional"
$ws.Cells.Item(161,4).NumberFormat = "@"
$ws.Cells.Item(161,4).Value = "2023"
$ws.Cells.Item(161,4).Style = "Normal"
$ws.Cells.Item(161,5).Value = 45235.85416666666
$ws.Cells.Item(161,6).Value = "Union de Santa Fe"
$ws.Cells.Item(161,7).Value = 0
$ws.Cells.Item(161,8).Value = "Lanus"
$ws.Cells.Item(161,9).Value = 0
$ws.Cells.Item(161,10).Value = 2.29
$ws.Cells.Item(161,11).Value = "01/11/2023 01:12"
$ws.Cells.Item(161,12).Value = 2.59
$ws.Cells.Item(161,13).Value = "05/11/2023 20:22"
$ws.Cells.Item(161,14).Value = 3.06
$ws.Cells.Item(161,15).Value = "01/11/2023 01:12"
$ws.Cells.Item(161,16).Value = 2.84
$ws.Cells.Item(161,17).Value = "05/11/2023 20:20"
$ws.Cells.Item(161,18).Value = 3.36
$ws.Cells.Item(161,19).Value = "01/11/2023 01:12"
$ws.Cells.Item(161,20).Value = 3.39
$ws.Cells.Item(161,21).Value = "05/11/2023 20:22"
$ws.Cells.Item(161,22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/union-de-santa-fe-lanus/hIyss4cF/"
